$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.165.91"
$ws.Range("E2").Value = "  -0.40%  "
$ws.Range("D3").Value = "2.021.80"
$ws.Range("E3").Value = "  -1.30%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.95"
$ws.Range("E5").Value = "  -1.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.609"
$ws.Range("E6").Value = "  -1.84%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "55.16"
$ws.Range("E8").Value = "  -3.29%  "
$ws.Range("E9").Value = "  -0.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0790"
$ws.Range("E10").Value = "  +2.81%  "
$ws.Range("D12").Value = "2.321.76"
$ws.Range("E12").Value = "  -1.19%  "
$ws.Range("E13").Value = "  -2.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.52"
$ws.Range("E14").Value = "  -0.17%  "
$ws.Range("E15").Value = "  -1.33%  "
$ws.Range("E16").Value = "  -1.69%  "
$ws.Range("D17").Value = "2.011.68"
$ws.Range("E17").Value = "  -1.89%  "
$ws.Range("D18").Value = "37.053.91"
$ws.Range("E18").Value = "  -0.68%  "
$ws.Range("E19").Value = "  +3.98%  "
$ws.Range("E20").Value = "  -0.94%  "
$ws.Range("D21").Value = "0.0₃0825"
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.31"
$ws.Range("E22").Value = "  +0.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("E24").Value = "  +2.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.21"
$ws.Range("E25").Value = "  -5.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.70"
$ws.Range("E26").Value = "  -1.75%  "
$ws.Range("E27").Value = "  -2.81%  "
$ws.Range("E28").Value = "  -2.31%  "
$ws.Range("E29").Value = "  +0.36%  "
$ws.Range("E30").Value = "  -2.05%  "
$ws.Range("E31").Value = "  -3.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.57"
$ws.Range("E32").Value = "  +0.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0618"
$ws.Range("E33").Value = "  -1.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.45"
$ws.Range("E34").Value = "  -2.57%  "
$ws.Range("E35").Value = "  -4.17%  "
$ws.Range("E36").Value = "  +1.75%  "
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.16"
$ws.Range("E38").Value = "  -3.75%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.43"
$ws.Range("E39").Value = "  +2.97%  "
$ws.Range("E40").Value = "  -3.63%  "
$ws.Range("D41").Value = "1.483.74"
$ws.Range("E41").Value = "  +0.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "16.68"
$ws.Range("E42").Value = "  +0.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "95.51"
$ws.Range("E43").Value = "  -2.45%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0929"
$ws.Range("E44").Value = "  -2.54%  "
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.14"
$ws.Range("E45").Value = "  -4.15%  "
$ws.Range("B46").Value = "HuobiToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.76"
$ws.Range("E46").Value = "  -4.86%  "
$ws.Range("E47").Value = "  +1.63%  "
$ws.Range("E48").Value = "  -1.78%  "
$ws.Range("E49").Value = "  -0.22%  "
$ws.Range("D50").Value = "2.208.92"
$ws.Range("E50").Value = "  -1.18%  "
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "44.55"
$ws.Range("E51").Value = "  -1.57%  "
